$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row 16: Martha Milena Villalba Borja record (moved to the top)
$ws.Range("C16").Value = "45747862"
$ws.Range("D16").Value = "MARTHA MILENA VILLALBA BORJA"
$ws.Range("E16").Value = "1709"
$ws.Range("F16").Value = 24591
$ws.Range("G16").Value = 737717

# Rows 17-25: Vladimir Genes Campo records, periods now ascending 2003-2011
$ws.Range("C17").Value = "3811636"
$ws.Range("D17").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E17").Value = "2003"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 877804

$ws.Range("C18").Value = "3811636"
$ws.Range("D18").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E18").Value = "2004"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 877804

$ws.Range("C19").Value = "3811636"
$ws.Range("D19").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E19").Value = "2005"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 877804

$ws.Range("C20").Value = "3811636"
$ws.Range("D20").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E20").Value = "2006"
$ws.Range("F20").Value = 35112
$ws.Range("G20").Value = 877804

$ws.Range("C21").Value = "3811636"
$ws.Range("D21").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E21").Value = "2007"
$ws.Range("F21").Value = 35112
$ws.Range("G21").Value = 877804

$ws.Range("C22").Value = "3811636"
$ws.Range("D22").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E22").Value = "2008"
$ws.Range("F22").Value = 35112
$ws.Range("G22").Value = 877804

$ws.Range("C23").Value = "3811636"
$ws.Range("D23").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E23").Value = "2009"
$ws.Range("F23").Value = 35112
$ws.Range("G23").Value = 877804

$ws.Range("C24").Value = "3811636"
$ws.Range("D24").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E24").Value = "2010"
$ws.Range("F24").Value = 35112
$ws.Range("G24").Value = 877804

$ws.Range("C25").Value = "3811636"
$ws.Range("D25").Value = "VLADIMIR GENES CAMPO"
$ws.Range("E25").Value = "2011"
$ws.Range("F25").Value = 35112
$ws.Range("G25").Value = 877804
